$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("BaseInformation").Range("A2:F2").ClearContents()
$wb.Worksheets.Item("Terms").Range("A2:G2").ClearContents()
$wb.Worksheets.Item("Labels").Range("A2:D2").ClearContents()
$wb.Worksheets.Item("References").Range("A2:I2").ClearContents()
$wb.Worksheets.Item("Structures").Range("A2:C2").ClearContents()
$wb.Worksheets.Item("Associations").Range("A2:H2").ClearContents()
$wb.Worksheets.Item("Rules-Consistency").Range("A2:G2").ClearContents()
$wb.Worksheets.Item("Rules-RollForward").Range("A2:G2").ClearContents()
$wb.Worksheets.Item("Rules-MemberAggregation").Range("A2:H2").ClearContents()
$wb.Worksheets.Item("Rules-Adjustment").Range("D2:J2").ClearContents()
$wb.Worksheets.Item("Rules-Variance").Range("D2:K2").ClearContents()
$wb.Worksheets.Item("Rules-Nonstandard").Range("A2:E2").ClearContents()
$wb.Worksheets.Item("Facts").Range("A2:H2").ClearContents()
$wb.Worksheets.Item("Facts-Dimensions").Range("A2:D2").ClearContents()

$factsParen = $wb.Worksheets.Item("Facts-Parenthetical")
$factsParen.Range("B2:D2").ClearContents()
$factsParen.Range("A2").Clear()
